$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Create new row 126 by duplicating row 125 (full row copy preserves styles/values)
$ws.Range("A125:R125").Copy($ws.Range("A126:R126"))

# Apply the weekly-shift update: row N gets the data previously in row N-1,
# and row 103 receives the new week's data.
$ws.Range("D103").Value = 45005
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 250
$ws.Range("K103").Value = 1300
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = 1400
$ws.Range("P103").Value = 700

$ws.Range("D104").Value = 44529
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 300
$ws.Range("K104").Value = 1800
$ws.Range("L104").Value = 2000
$ws.Range("M104").Value = 1900
$ws.Range("P104").Value = 950

$ws.Range("D105").Value = 44901
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 350
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 3500
$ws.Range("M105").Value = 3214
$ws.Range("P105").Value = 1607

$ws.Range("D106").Value = 44859
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 250
$ws.Range("K106").Value = 1400
$ws.Range("L106").Value = 1500
$ws.Range("M106").Value = 1450
$ws.Range("P106").Value = 725

$ws.Range("D107").Value = 44169
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 2500
$ws.Range("M107").Value = 2250
$ws.Range("P107").Value = 1125

$ws.Range("D108").Value = 44165
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 300
$ws.Range("K108").Value = 1000
$ws.Range("L108").Value = 1200
$ws.Range("M108").Value = 1100
$ws.Range("P108").Value = 550

$ws.Range("D109").Value = 44572
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 300
$ws.Range("K109").Value = 4000
$ws.Range("L109").Value = 4500
$ws.Range("M109").Value = 4250
$ws.Range("P109").Value = 2125

$ws.Range("D110").Value = 44760
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 300
$ws.Range("K110").Value = 3000
$ws.Range("L110").Value = 3500
$ws.Range("M110").Value = 3250
$ws.Range("P110").Value = 1625

$ws.Range("D111").Value = 44876
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 200
$ws.Range("K111").Value = 1000
$ws.Range("L111").Value = 1200
$ws.Range("M111").Value = 1100
$ws.Range("P111").Value = 550

$ws.Range("D112").Value = 44386
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 250
$ws.Range("K112").Value = 3500
$ws.Range("L112").Value = 4000
$ws.Range("M112").Value = 3750
$ws.Range("P112").Value = 1875

$ws.Range("D113").Value = 44428
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 270
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 3800
$ws.Range("M113").Value = 3650
$ws.Range("P113").Value = 1825

$ws.Range("D114").Value = 44253
$ws.Range("I114").Value = "Primera"
$ws.Range("J114").Value = 300
$ws.Range("K114").Value = 2400
$ws.Range("L114").Value = 2500
$ws.Range("M114").Value = 2450
$ws.Range("P114").Value = 1225

$ws.Range("D115").Value = 44998
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 270
$ws.Range("K115").Value = 1500
$ws.Range("L115").Value = 2000
$ws.Range("M115").Value = 1750
$ws.Range("P115").Value = 875

$ws.Range("D116").Value = 44740
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 300
$ws.Range("K116").Value = 3500
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = 3750
$ws.Range("P116").Value = 1875

$ws.Range("D117").Value = 44176
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 300
$ws.Range("K117").Value = 1900
$ws.Range("L117").Value = 2000
$ws.Range("M117").Value = 1950
$ws.Range("P117").Value = 975

$ws.Range("D118").Value = 44795
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 300
$ws.Range("K118").Value = 2000
$ws.Range("L118").Value = 2500
$ws.Range("M118").Value = 2250
$ws.Range("P118").Value = 1125

$ws.Range("D119").Value = 44673
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 250
$ws.Range("K119").Value = 2000
$ws.Range("L119").Value = 2500
$ws.Range("M119").Value = 2250
$ws.Range("P119").Value = 1125

$ws.Range("D120").Value = 44568
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 300
$ws.Range("K120").Value = 5000
$ws.Range("L120").Value = 5500
$ws.Range("M120").Value = 5250
$ws.Range("P120").Value = 2625

$ws.Range("D121").Value = 44473
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 300
$ws.Range("K121").Value = 950
$ws.Range("L121").Value = 1000
$ws.Range("M121").Value = 975
$ws.Range("P121").Value = 488

$ws.Range("D122").Value = 44606
$ws.Range("I122").Value = "Segunda"
$ws.Range("J122").Value = 200
$ws.Range("K122").Value = 2000
$ws.Range("L122").Value = 2500
$ws.Range("M122").Value = 2250
$ws.Range("P122").Value = 1125

$ws.Range("D123").Value = 44677
$ws.Range("I123").Value = "Primera"
$ws.Range("J123").Value = 300
$ws.Range("K123").Value = 2300
$ws.Range("L123").Value = 2500
$ws.Range("M123").Value = 2400
$ws.Range("P123").Value = 1200

$ws.Range("D124").Value = 44747
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 300
$ws.Range("K124").Value = 3500
$ws.Range("L124").Value = 4000
$ws.Range("M124").Value = 3750
$ws.Range("P124").Value = 1875

$ws.Range("D125").Value = 44498
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 300
$ws.Range("K125").Value = 800
$ws.Range("L125").Value = 1000
$ws.Range("M125").Value = 900
$ws.Range("P125").Value = 450

$ws.Range("D126").Value = 44608
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 300
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 3500
$ws.Range("M126").Value = 3250
$ws.Range("P126").Value = 1625
